$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New table contents (header + data rows) replacing the old ones.
$data = @(
    @("Name", "Type"),
    @("participant_id", "integer"),
    @("platekey", "character varying"),
    @("referral_id", "character varying"),
    @("associated_interpretation_request_id", "character varying"),
    @("delivery_type", "character varying"),
    @("delivery_id", "character varying"),
    @("delivery_date", "timestamp without time zone"),
    @("path", "character varying"),
    @("delivery_version", "character varying"),
    @("genome_build", "character varying"),
    @("data_format", "character varying")
)

$rowCount = $data.Length

# Copy the existing data-row formatting down onto the new last row before
# writing values, so the appended row matches the rest of the table.
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)

# Write order matches how the strings first appear in the shared string
# table on the source workbook: rows 1-5 introduce their new text first,
# then row 12's "data_format" is registered, then row 6's "delivery_type",
# and finally the remaining (already-known) rows fill in.
$order = @(0, 1, 2, 3, 4, 11, 5, 6, 7, 8, 9, 10)

foreach ($i in $order) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Column widths to fit the new (longer) text (closest achievable values;
# the host quantizes ColumnWidth to 1/6-character steps).
$ws.Columns.Item(1).ColumnWidth = 23.983072916666668
$ws.Columns.Item(2).ColumnWidth = 25.436197916666668

# Match the new active selection cell.
$ws.Range("F7").Select()
